# Update the legacy GSC export (gsc-export-old/HTTPS.xlsx): the rolling
# date window advanced by one day.
#   - "Chart" sheet: drop the oldest date (2025-09-18), shift every
#     Non-HTTPS/HTTPS daily count up by one row, and append a new row for
#     the newly-entered day (2025-12-15) with zero counts.
#   - "Table" sheet header is untouched; its shared-string indices just
#     follow the shared-string-table shift automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# New column C ("HTTPS URLs") values for rows 2..89: the old column C
# (rows 3..88) shifted up by one row, with two trailing zeros appended
# (the carried-forward last value, plus the brand new day).
$newC = @(24.0,34.0,34.0,34.0,33.0,44.0,38.0,44.0,39.0,34.0,46.0,46.0,56.0,67.0,67.0,74.0,83.0,78.0,78.0,73.0,70.0,65.0,60.0,57.0,50.0,47.0,39.0,36.0,26.0,23.0,15.0,5.0,5.0,5.0,4.0,3.0,2.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0,0.0)

# Date column starts the day after the old first row (2025-09-18 is
# dropped) and runs through the two new trailing days.
$date = Get-Date -Year 2025 -Month 9 -Day 19

for ($i = 0; $i -lt $newC.Length; $i++) {
    $r = $i + 2
    $dateText = "'" + $date.ToString("yyyy-MM-dd")
    $ws.Cells.Item($r, 1).Value = $dateText
    $ws.Cells.Item($r, 2).Value = 0.0
    $ws.Cells.Item($r, 3).Value = $newC[$i]
    $date = $date.AddDays(1)
}

# The freshly-typed dates land as quote-prefixed text with a style
# override (Excel always marks a leading-apostrophe text entry that
# way); strip that back off so column A keeps the sheet's plain default
# formatting, same as every other cell.
$ws.Range("A2:A89").ClearFormats()
